$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# The "преподователь" (teacher) column used the placeholder value "n" for
# several lessons that still needed an assigned teacher. Rename that
# placeholder in place (this updates every cell still showing "n", e.g.
# D2, D3, D18, D26) to the first real name that fills in one of those
# slots.
$ws.Cells.Replace("n", "Шапошникова И.В")

# Now assign the remaining, still-unresolved lessons their own specific
# teacher. Order matters only for matching shared-string insertion order,
# first-use order below mirrors the original authoring order.
$ws.Range("D6").Value  = "Семенов О.Ю"
$ws.Range("D11").Value = "Кулагина И.В"
$ws.Range("D9").Value  = "Пешков А.А"
$ws.Range("D30").Value = "Костюнина М.В"
$ws.Range("D5").Value  = "Кузнецова С.В"

$ws.Range("D10").Value = "Семенов О.Ю"
$ws.Range("D19").Value = "Семенов О.Ю"
$ws.Range("D21").Value = "Семенов О.Ю"
$ws.Range("D22").Value = "Кулагина И.В"
$ws.Range("D23").Value = "Пешков А.А"
$ws.Range("D33").Value = "Пешков А.А"

$ws.Range("D5").Select()
